# Updated cryptos list on Fri Jun  9 08:10:33 UTC 2023 with GitHub Actions
#
# Note: the "Price"/"Volume(1h)" columns are stored as plain text in this
# workbook (General number format). Excel's Range.Value setter auto-detects
# numeric-looking strings and silently converts them to real numbers, which
# would corrupt values like "18.60" -> 18.6 or "0.07750" -> 0.0775. Prefixing
# the literal with a leading apostrophe forces Excel to keep it as text,
# exactly like typing it in by hand, while leaving the General number format
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 35 / 36: HuobiToken and ImmutableX swap places, values refresh ---
# Before: row35 = HuobiToken / row36 = ImmutableX
# After:  row35 = ImmutableX / row36 = HuobiToken
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7169"
$ws.Range("E35").Value = "  +0.48%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.856"
$ws.Range("E36").Value = "  +0.49%  "

# --- Row 48 / 49: Cronos and Algorand swap places, values refresh ---
# Before: row48 = Cronos / row49 = Algorand
# After:  row48 = Algorand / row49 = Cronos
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1226"
$ws.Range("E48").Value = "  +0.82%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05805"
$ws.Range("E49").Value = "  -1.19%  "

# --- Remaining price / volume refreshes ---
$ws.Range("D2").Value = "26.516.08"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.836.83"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'258.31"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.5251"
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("D8").Value = "'0.3141"
$ws.Range("E8").Value = "  -4.00%  "
$ws.Range("D9").Value = "'0.06777"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "'18.60"
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "'0.7748"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "'0.07750"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").Value = "1.832.32"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "'87.54"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").Value = "'4.995"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("D17").Value = "'13.78"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "'0.000007901"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "26.530.62"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "2.067.76"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "'4.584"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'5.939"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("D24").Value = "'9.295"
$ws.Range("E24").Value = "  -1.99%  "
$ws.Range("D25").Value = "'142.53"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("D26").Value = "'2.190"
$ws.Range("D27").Value = "'1.680"
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("D29").Value = "'110.43"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").Value = "'4.151"
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("D31").Value = "'0.08705"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "'4.042"
$ws.Range("E32").Value = "  -2.10%  "
$ws.Range("D33").Value = "'0.04849"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("D37").Value = "'3.081"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'2.223"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "'0.01717"
$ws.Range("E39").Value = "  -3.03%  "
$ws.Range("D40").Value = "'0.4793"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").Value = "'0.8907"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("D42").Value = "'109.61"
$ws.Range("E42").Value = "  -2.01%  "
$ws.Range("D43").Value = "'5.903"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'7.598"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").Value = "'0.4137"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").Value = "'8.956"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").Value = "'0.8919"
$ws.Range("E51").Value = "  +0.38%  "
